# "Add menu button in the game"
# Adds a new dialog row (dialog_back_to_menu) to the defaultDialog sheet,
# mirroring the other rows: dialogId text in col A (shared string), then
# four numeric/zero columns with the npcType flag of 4 in column B.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRow = 48

$ws.Range("A" + $newRow).Value = "dialog_back_to_menu"
$ws.Range("B" + $newRow).Value = 4
$ws.Range("C" + $newRow).Value = 0
$ws.Range("D" + $newRow).Value = 0
$ws.Range("E" + $newRow).Value = 0
$ws.Range("F" + $newRow).Value = 0

# Match the author's final on-screen selection/scroll state after adding
# the row (select the newly populated data cells of the new row).
$ws.Range("B48:F48").Select() | Out-Null
